$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the BSURL column (I) entirely - it is no longer part of the
# environment data table.
$ws.Columns("I").Delete()

# Point the environment URLs at the new "sandbox" environment instead of
# the old "test19" one.
$ws.Range("A2").Value = "https://sandbox.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://sandbox.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://sandbox.cliotest.com/warehouse/control/main"
$ws.Range("E2").Value = "https://mirandakate.cabisandbox.com"
$ws.Range("F2").Value = "virtual_cabitest21"
$ws.Range("G2").Value = "sandbox"
$ws.Range("H2").Value = "cabisandbox"

# Rebuild the hyperlinks (same targets as before) so they end up attached
# to the refreshed cells in natural column order.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://sandbox.cabiclio.com/backoffice/control/main")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://test17.cliotest.com/cabicentral/control/main", "", "", "https://test17.cliotest.com/cabicentral/control/main")
$ws.Hyperlinks.Add($ws.Range("D2"), "https://sandbox.cabiclio.com/warehouse/control/main", "", "", "https://sandbox.cabiclio.com/warehouse/control/main")
$ws.Hyperlinks.Add($ws.Range("E2"), "https://mirandakate.cabitest5.com/")

# Update the selection/view state to match the saved workbook.
$ws.Range("H2").Select()
$excel.ActiveWindow.ScrollColumn = 3
